$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.759.88'
$ws.Range("E2").Value = '  -0.29%  '
$ws.Range("D3").Value = '1.867.23'
$ws.Range("E3").Value = '  -0.11%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9988'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7292'
$ws.Range("E5").Value = '  -1.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '240.41'
$ws.Range("E6").Value = '  -0.61%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9997'
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3118'
$ws.Range("E8").Value = '  -1.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07082'
$ws.Range("E9").Value = '  -0.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.26'
$ws.Range("E10").Value = '  -1.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08199'
$ws.Range("E11").Value = '  -2.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7395'
$ws.Range("E12").Value = '  -1.60%  '
$ws.Range("D13").Value = '1.877.48'
$ws.Range("E13").Value = '  +0.58%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.304'
$ws.Range("E14").Value = '  -2.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.93'
$ws.Range("E15").Value = '  -0.55%  '
$ws.Range("D16").Value = '29.787.78'
$ws.Range("E16").Value = '  -0.23%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.997'
$ws.Range("E17").Value = '  -0.44%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '247.40'
$ws.Range("E18").Value = '  +1.97%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.34'
$ws.Range("E19").Value = '  -1.63%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007789'
$ws.Range("E20").Value = '  -0.32%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.153.06'
$ws.Range("E21").Value = '  +1.78%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.22%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9988'
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.725'
$ws.Range("E24").Value = '  -2.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1535'
$ws.Range("E25").Value = '  -1.64%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.155'
$ws.Range("E26").Value = '  -1.44%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '163.12'
$ws.Range("E27").Value = '  -0.60%  '
$ws.Range("E28").Value = '  -0.50%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.001'
$ws.Range("E29").Value = '  -0.74%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.439'
$ws.Range("E30").Value = '  -2.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.511'
$ws.Range("E31").Value = '  -2.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.517'
$ws.Range("E32").Value = '  -1.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.158'
$ws.Range("E33").Value = '  -3.16%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05261'
$ws.Range("E34").Value = '  -0.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.229'
$ws.Range("E35").Value = '  -0.24%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7443'
$ws.Range("E36").Value = '  -1.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.001'
$ws.Range("E37").Value = '  +0.11%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.687'
$ws.Range("E38").Value = '  -0.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01927'
$ws.Range("E39").Value = '  -1.41%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.728'
$ws.Range("E40").Value = '  -0.79%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4445'
$ws.Range("E41").Value = '  -0.55%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.983'
$ws.Range("E42").Value = '  -1.42%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8650'
$ws.Range("E43").Value = '  +0.66%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '71.02'
$ws.Range("E44").Value = '  -1.54%  '
$ws.Range("D45").Value = '1.040.39'
$ws.Range("E45").Value = '  -5.52%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9994'
$ws.Range("E46").Value = '  -0.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '103.78'
$ws.Range("E47").Value = '  +0.66%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.813'
$ws.Range("E48").Value = '  -1.41%  '
$ws.Range("B49").Value = 'Aptos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.442'
$ws.Range("E49").Value = '  -3.03%  '
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.048.82'
$ws.Range("E50").Value = '  +1.70%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.456'
$ws.Range("E51").Value = '  -0.58%  '
